$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AR")
$ws1.Range("B2").Value = 0.0008887854659561153
$ws1.Range("B3").Value = 0.7727922900322718
$ws1.Range("B4").Value = 0.1005288770640802

$ws2 = $wb.Worksheets.Item("SETAR")
$ws2.Range("B2").Value = -0.09306512800291195
$ws2.Range("B3").Value = 0.729940048388102
$ws2.Range("B4").Value = 0.07780760966866578
$ws2.Range("B5").Value = 0.09573379536244561
$ws2.Range("B6").Value = 0.6962706411279245
$ws2.Range("B7").Value = 0.06103817075380358

$ws3 = $wb.Worksheets.Item("GARCH")
$ws3.Range("B2").Value = -0.0009891166833161688
$ws3.Range("B3").Value = 0.002268313855635574
$ws3.Range("B4").Value = 0.1999993330331042
$ws3.Range("B5").Value = 0.7800007428514291

$ws4 = $wb.Worksheets.Item("TARCH")
$ws4.Range("B2").Value = -0.001065093801448775
$ws4.Range("B3").Value = 0.002269123102064011
$ws4.Range("B4").Value = 0.1999995091724977
$ws4.Range("B5").Value = 0.00999997185428371
$ws4.Range("B6").Value = 0.7750005422751267

$ws5 = $wb.Worksheets.Item("AR-TARCH")
$ws5.Range("B2").Value = 0.001706453550795025
$ws5.Range("B3").Value = 0.7760529830641222
$ws5.Range("B4").Value = 0.002010370423522397
$ws5.Range("B5").Value = 0.199979546336019
$ws5.Range("B6").Value = 0.009999039529196617
$ws5.Range("B7").Value = 0.775023178175507
